$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.15
$ws.Range("H2").Value = 2.65
$ws.Range("O2").Value = 1.33
$ws.Range("G3").Value = 1.82
$ws.Range("I3").Value = 4.6
$ws.Range("L3").Value = 1.52
$ws.Range("M3").Value = 2.22
$ws.Range("N3").Value = 2.47
$ws.Range("Q3").Value = 2.1
$ws.Range("T3").Value = 4.85
$ws.Range("X3").Value = 19
$ws.Range("Z3").Value = 5.1
$ws.Range("AB3").Value = 23
$ws.Range("AC3").Value = 175
$ws.Range("AD3").Value = 9
$ws.Range("AE3").Value = 24
$ws.Range("AF3").Value = 17
$ws.Range("AH3").Value = 65
$ws.Range("AI3").Value = 90
$ws.Range("H4").Value = 3.05
$ws.Range("I4").Value = 2.87
$ws.Range("L4").Value = 1.42
$ws.Range("M4").Value = 2.45
$ws.Range("N4").Value = 2.22
$ws.Range("O4").Value = 1.52
$ws.Range("P4").Value = 1.5
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 1.93
$ws.Range("S4").Value = 1.7
$ws.Range("T4").Value = 6.6
$ws.Range("X4").Value = 23
$ws.Range("Y4").Value = 40
$ws.Range("Z4").Value = 7.3
$ws.Range("AA4").Value = 6
$ws.Range("AB4").Value = 17
$ws.Range("AC4").Value = 100
$ws.Range("AD4").Value = 7.2
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 11
$ws.Range("AH4").Value = 29
$ws.Range("AI4").Value = 45
$ws.Range("AJ4").Value = 800
$ws.Range("G6").Value = 1.5
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 1.05
$ws.Range("K6").Value = 11
$ws.Range("R6").Value = 2.1
$ws.Range("S6").Value = 1.63
$ws.Range("W6").Value = 10
$ws.Range("Y6").Value = 29
$ws.Range("AA6").Value = 8.5
$ws.Range("AD6").Value = 12
$ws.Range("AF6").Value = 17
$ws.Range("AH6").Value = 41
$ws.Range("P7").Value = 1.37
$ws.Range("G11").Value = 1.72
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 5.3
$ws.Range("J11").Value = 1.1
$ws.Range("K11").Value = 5.9
$ws.Range("L11").Value = 1.44
$ws.Range("M11").Value = 2.6
$ws.Range("N11").Value = 2.3
$ws.Range("O11").Value = 1.55
$ws.Range("P11").Value = 1.47
$ws.Range("Q11").Value = 2.5
$ws.Range("R11").Value = 2.12
$ws.Range("S11").Value = 1.65
$ws.Range("U11").Value = 6.9
$ws.Range("V11").Value = 8.5
$ws.Range("W11").Value = 13.5
$ws.Range("X11").Value = 16.5
$ws.Range("Y11").Value = 37
$ws.Range("Z11").Value = 5.9
$ws.Range("AA11").Value = 6.4
$ws.Range("AB11").Value = 19
$ws.Range("AC11").Value = 120
$ws.Range("AD11").Value = 11.25
$ws.Range("AE11").Value = 29
$ws.Range("AF11").Value = 17
$ws.Range("AG11").Value = 110
$ws.Range("AH11").Value = 65
$ws.Range("AI11").Value = 70
$ws.Range("G12").Value = 4.6
$ws.Range("H12").Value = 3.5
$ws.Range("I12").Value = 1.72
$ws.Range("J12").Value = 1.06
$ws.Range("K12").Value = 7.3
$ws.Range("L12").Value = 1.29
$ws.Range("M12").Value = 3.25
$ws.Range("N12").Value = 1.87
$ws.Range("O12").Value = 1.85
$ws.Range("P12").Value = 1.44
$ws.Range("Q12").Value = 2.62
$ws.Range("R12").Value = 1.8
$ws.Range("S12").Value = 1.91
$ws.Range("T12").Value = 12
$ws.Range("U12").Value = 26
$ws.Range("V12").Value = 15
$ws.Range("W12").Value = 80
$ws.Range("X12").Value = 45
$ws.Range("Y12").Value = 50
$ws.Range("Z12").Value = 7.3
$ws.Range("AA12").Value = 6.8
$ws.Range("AB12").Value = 15
$ws.Range("AC12").Value = 70
$ws.Range("AE12").Value = 8.25
$ws.Range("AF12").Value = 8
$ws.Range("AG12").Value = 14
$ws.Range("AH12").Value = 13.5
$ws.Range("AI12").Value = 25
$ws.Range("AJ12").Value = 600
$ws.Range("G13").Value = 3.6
$ws.Range("H13").Value = 2.75
$ws.Range("I13").Value = 2.3
$ws.Range("L13").Value = 1.34
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 2
$ws.Range("O13").Value = 1.72
$ws.Range("P13").Value = 1.38
$ws.Range("Q13").Value = 2.8
$ws.Range("R13").Value = 1.65
$ws.Range("S13").Value = 2.1
$ws.Range("T13").Value = 10.25
$ws.Range("U13").Value = 21
$ws.Range("V13").Value = 11.25
$ws.Range("W13").Value = 60
$ws.Range("AA13").Value = 5.3
$ws.Range("AB13").Value = 11.5
$ws.Range("AC13").Value = 50
$ws.Range("AD13").Value = 7.6
$ws.Range("AE13").Value = 11.75
$ws.Range("AG13").Value = 26
$ws.Range("AH13").Value = 18.5
$ws.Range("AI13").Value = 26
$ws.Range("AJ13").Value = 400
$ws.Range("G14").Value = 2.25
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 2.8
$ws.Range("N14").Value = 1.65
$ws.Range("O14").Value = 1.98
$ws.Range("T14").Value = 8
$ws.Range("U14").Value = 10.25
$ws.Range("W14").Value = 19
$ws.Range("X14").Value = 14
$ws.Range("Y14").Value = 18.5
$ws.Range("Z14").Value = 11.75
$ws.Range("AA14").Value = 5.8
$ws.Range("AB14").Value = 9.75
$ws.Range("AC14").Value = 32
$ws.Range("AD14").Value = 9.5
$ws.Range("AF14").Value = 8.5
$ws.Range("AG14").Value = 28
$ws.Range("AI14").Value = 20
$ws.Range("AJ14").Value = 175
$ws.Range("AA15").Value = 9.75
$ws.Range("AB15").Value = 25
$ws.Range("AC15").Value = 120
$ws.Range("AD15").Value = 26
$ws.Range("AE15").Value = 90
$ws.Range("AI15").Value = 150
$ws.Range("G16").Value = 1.7
$ws.Range("H16").Value = 3.3
$ws.Range("I16").Value = 4.65
$ws.Range("N16").Value = 2.05
$ws.Range("O16").Value = 1.6
$ws.Range("P16").Value = 1.39
$ws.Range("Q16").Value = 2.42
$ws.Range("T16").Value = 4.9
$ws.Range("U16").Value = 6
$ws.Range("V16").Value = 7.1
$ws.Range("W16").Value = 10.25
$ws.Range("X16").Value = 12.5
$ws.Range("Y16").Value = 26
$ws.Range("Z16").Value = 7.9
$ws.Range("AA16").Value = 5.8
$ws.Range("AB16").Value = 14.5
$ws.Range("AC16").Value = 75
$ws.Range("AD16").Value = 9.25
$ws.Range("AE16").Value = 20
$ws.Range("AF16").Value = 13
$ws.Range("AG16").Value = 60
$ws.Range("AH16").Value = 40
$ws.Range("AI16").Value = 45
$ws.Range("AJ16").Value = 500
$ws.Range("K17").Value = 6.5
$ws.Range("L17").Value = 1.36
$ws.Range("M17").Value = 2.9
$ws.Range("N17").Value = 2.05
$ws.Range("R17").Value = 1.83
$ws.Range("S17").Value = 1.87
$ws.Range("Y17").Value = 30
$ws.Range("Z17").Value = 6.5
$ws.Range("AB17").Value = 15
$ws.Range("AC17").Value = 75
$ws.Range("AD17").Value = 9.75
$ws.Range("AE17").Value = 19.5
$ws.Range("AG17").Value = 50
$ws.Range("AH17").Value = 37
$ws.Range("AI17").Value = 45
$ws.Range("AJ17").Value = 700
$ws.Range("H18").Value = 2.87
$ws.Range("G24").Value = 2.8
$ws.Range("I24").Value = 2.55
$ws.Range("K24").Value = 8.5
$ws.Range("T24").Value = 8
$ws.Range("X24").Value = 23
$ws.Range("AE24").Value = 12
$ws.Range("AG24").Value = 26
$ws.Range("G25").Value = 1.49
$ws.Range("R25").Value = 1.41
$ws.Range("S25").Value = 2.62
$ws.Range("J28").Value = 1.11
$ws.Range("K29").Value = 6
$ws.Range("N29").Value = 2.32
$ws.Range("Z29").Value = 6
$ws.Range("G31").Value = 2.25
$ws.Range("H31").Value = 3.25
$ws.Range("I31").Value = 3.2
$ws.Range("N31").Value = 2.1
$ws.Range("O31").Value = 1.7
$ws.Range("T31").Value = 7.5
$ws.Range("U31").Value = 10
$ws.Range("V31").Value = 9.5
$ws.Range("W31").Value = 21
$ws.Range("X31").Value = 19
$ws.Range("Y31").Value = 29
$ws.Range("AD31").Value = 9
$ws.Range("AE31").Value = 15
$ws.Range("AF31").Value = 12
$ws.Range("AG31").Value = 34
$ws.Range("AH31").Value = 26
$ws.Range("I34").Value = 2.63
$ws.Range("J34").Value = 1.02
$ws.Range("L34").Value = 1.15
$ws.Range("T34").Value = 11
$ws.Range("P35").Value = 1.33
$ws.Range("J36").Value = 1.04
$ws.Range("K36").Value = 13
$ws.Range("AD36").Value = 17
$ws.Range("AI36").Value = 41
$ws.Range("T38").Value = 8
$ws.Range("W38").Value = 15
$ws.Range("AH38").Value = 34
$ws.Range("K40").Value = 13
$ws.Range("L40").Value = 1.22
$ws.Range("M40").Value = 4
$ws.Range("R40").Value = 1.62
$ws.Range("S40").Value = 2.2
$ws.Range("AB40").Value = 12
$ws.Range("G42").Value = 1.95
$ws.Range("I42").Value = 3.9
$ws.Range("J42").Value = 1.06
$ws.Range("K42").Value = 10
$ws.Range("R42").Value = 1.75
$ws.Range("S42").Value = 2
$ws.Range("U42").Value = 9.5
$ws.Range("P43").Value = 1.3
$ws.Range("N44").Value = 1.44
$ws.Range("G49").Value = 2.37
$ws.Range("L49").Value = 1.31
$ws.Range("Q49").Value = 2.82
$ws.Range("X49").Value = 20
$ws.Range("AF49").Value = 10
